# Estadisticos Matutinos 15 Oct
# Updates the grade-progress stats on "Estadisticos 1P" / "Estadisticos Final"
# (row 2, the NC / 3BEM group), adds the matching partial-count on
# "Estadisticos 2P", and appends one rescued-student row on "Rescatables".

$wb = $excel.ActiveWorkbook

# --- Estadisticos 1P: row 2 (NC / 3BEM) -----------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 7
$ws1.Range("F2").Value = 25
$ws1.Range("G2").Value = 78.13
$ws1.Range("H2").Value = 8.199999999999999

# --- Estadisticos 2P: row 2 (NC / 3BEM) -----------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 25

# --- Estadisticos Final: row 2 (NC / 3BEM) --------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 7
$ws3.Range("F2").Value = 25
$ws3.Range("G2").Value = 78.13
$ws3.Range("H2").Value = 8.199999999999999

# --- Rescatables: append the new rescued student --------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Range("A2").Value = 19330051920085
$ws4.Range("B2").Value = "ANTONIO"
$ws4.Range("C2").Value = "TEXOCO"
$ws4.Range("D2").Value = "JOSE JAZAEL"
$ws4.Range("E2").Value = "DISEÑA Y MANTIENE LOS SISTEMAS DE ILUMINACIÓN"
$ws4.Range("F2").Value = "5BEM"
$ws4.Range("G2").Value = 7
